$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three "Шаги воспроизведения" (steps) cells: "Eg" -> "Eng" (and "Ru" -> "Рус")
$ws.Range("E3").Value = "1. Выбрать `"Eng`" `n2. Выбрать слайд (class=`"swiper-wrapper`") номер 6"
$ws.Range("E4").Value = "1. Выбрать `"Eng`" "
$ws.Range("E5").Value = "1. Выбрать `"Eng`" `n2. Наблюдать кнопку `"Access for buyers`" с иконкой `"ключ`"`n3. Выбрать `"Рус`""

# Row heights changed (rows 3 and 4 shrank from 72 to 57.6)
$ws.Rows.Item(3).RowHeight = 57.6
$ws.Rows.Item(4).RowHeight = 57.6

# Update the active selection to E6
$ws.Range("E6").Select() | Out-Null
